$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 4,23
$data[0,0] = 0.00969425801640567
$data[0,1] = 0.0126771066368382
$data[0,2] = 0.938105891126025
$data[0,3] = 0.0126771066368382
$data[0,4] = 0.979120059656972
$data[0,5] = 0.968680089485459
$data[0,6] = 0.0111856823266219
$data[0,7] = 0.00969425801640567
$data[0,8] = 0.00149142431021626
$data[0,9] = 0.991051454138702
$data[0,10] = 0
$data[0,11] = 0.00447427293064877
$data[0,12] = 0
$data[0,13] = 0
$data[0,14] = 0
$data[0,15] = 0
$data[0,16] = 0.0380313199105145
$data[0,17] = 0.029082774049217
$data[0,18] = 0.00149142431021626
$data[0,19] = 0.0134228187919463
$data[0,20] = 0.0149142431021626
$data[0,21] = 0.00969425801640567
$data[0,22] = 0.00894854586129754
$data[1,0] = 0.967188665175242
$data[1,1] = 0.00671140939597315
$data[1,2] = 0.0104399701715138
$data[1,3] = 0.96793437733035
$data[1,4] = 0.0134228187919463
$data[1,5] = 0.00223713646532438
$data[1,6] = 0
$data[1,7] = 0.000745712155108128
$data[1,8] = 0
$data[1,9] = 0.0052199850857569
$data[1,10] = 0.994034302759135
$data[1,11] = 0.00447427293064877
$data[1,12] = 0.00223713646532438
$data[1,13] = 0.00149142431021626
$data[1,14] = 0.998508575689784
$data[1,15] = 0.991051454138702
$data[1,16] = 0.00447427293064877
$data[1,17] = 0.953020134228188
$data[1,18] = 0.0186428038777032
$data[1,19] = 0.000745712155108128
$data[1,20] = 0
$data[1,21] = 0.000745712155108128
$data[1,22] = 0.00372856077554064
$data[2,0] = 0.00298284862043251
$data[2,1] = 0.00969425801640567
$data[2,2] = 0.0454884414615958
$data[2,3] = 0.0052199850857569
$data[2,4] = 0.00223713646532438
$data[2,5] = 0.0260999254287845
$data[2,6] = 0.98806860551827
$data[2,7] = 0.988814317673378
$data[2,8] = 0.997762863534676
$data[2,9] = 0.00298284862043251
$data[2,10] = 0
$data[2,11] = 0.000745712155108128
$data[2,12] = 0.000745712155108128
$data[2,13] = 0.000745712155108128
$data[2,14] = 0
$data[2,15] = 0.00447427293064877
$data[2,16] = 0.957494407158837
$data[2,17] = 0.00223713646532438
$data[2,18] = 0.00223713646532438
$data[2,19] = 0.982102908277405
$data[2,20] = 0.979865771812081
$data[2,21] = 0.98806860551827
$data[2,22] = 0.985831469052946
$data[3,0] = 0.0201342281879195
$data[3,1] = 0.970917225950783
$data[3,2] = 0.00596569724086503
$data[3,3] = 0.0141685309470544
$data[3,4] = 0.0052199850857569
$data[3,5] = 0.00298284862043251
$data[3,6] = 0
$data[3,7] = 0.000745712155108128
$data[3,8] = 0.000745712155108128
$data[3,9] = 0.000745712155108128
$data[3,10] = 0.0052199850857569
$data[3,11] = 0.990305741983594
$data[3,12] = 0.997017151379567
$data[3,13] = 0.997762863534676
$data[3,14] = 0.00149142431021626
$data[3,15] = 0.00447427293064877
$data[3,16] = 0
$data[3,17] = 0.0156599552572707
$data[3,18] = 0.977628635346756
$data[3,19] = 0.00298284862043251
$data[3,20] = 0.00447427293064877
$data[3,21] = 0.000745712155108128
$data[3,22] = 0.00149142431021626

$ws.Range("B2:X5").Value = $data

